$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44159
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7455
$ws.Range("N2").Value = '$/caja 36 atados'
$ws.Range("O2").Value = 'Región Metropolitana'
$ws.Range("P2").Value = 207
$ws.Range("Q2").Value = 36

$ws.Range("D3").Value = 44263
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1914
$ws.Range("N3").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O3").Value = 'Provincia de Diguillín'
$ws.Range("P3").Value = 1914
$ws.Range("Q3").Value = 1

$ws.Range("D4").Value = 44272
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 1800
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1893
$ws.Range("N4").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O4").Value = 'Provincia de Diguillín'
$ws.Range("P4").Value = 1893
$ws.Range("Q4").Value = 1

$ws.Range("D5").Value = 44271
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1920
$ws.Range("N5").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O5").Value = 'Provincia de Diguillín'
$ws.Range("P5").Value = 1920
$ws.Range("Q5").Value = 1

$ws.Range("D6").Value = 44208
$ws.Range("J6").Value = 130
$ws.Range("K6").Value = 1800
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1908
$ws.Range("N6").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O6").Value = 'Provincia de Cautín'
$ws.Range("P6").Value = 1908
$ws.Range("Q6").Value = 1

$ws.Range("D7").Value = 44266
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 1913
$ws.Range("N7").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O7").Value = 'Provincia de Diguillín'
$ws.Range("P7").Value = 1913
$ws.Range("Q7").Value = 1

$ws.Range("D8").Value = 44265
$ws.Range("J8").Value = 220
$ws.Range("K8").Value = 1800
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 1909
$ws.Range("N8").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O8").Value = 'Provincia de Diguillín'
$ws.Range("P8").Value = 1909
$ws.Range("Q8").Value = 1

$ws.Range("D9").Value = 44270
$ws.Range("J9").Value = 260
$ws.Range("K9").Value = 1800
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = 1908
$ws.Range("N9").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O9").Value = 'Provincia de Diguillín'
$ws.Range("P9").Value = 1908
$ws.Range("Q9").Value = 1

$ws.Range("D10").Value = 44260
$ws.Range("J10").Value = 220
$ws.Range("K10").Value = 1800
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 1909
$ws.Range("N10").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O10").Value = 'Provincia de Diguillín'
$ws.Range("P10").Value = 1909
$ws.Range("Q10").Value = 1

$ws.Range("D11").Value = 44267
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1913
$ws.Range("N11").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O11").Value = 'Provincia de Diguillín'
$ws.Range("P11").Value = 1913
$ws.Range("Q11").Value = 1

$ws.Range("D12").Value = 44264
$ws.Range("J12").Value = 130
$ws.Range("K12").Value = 1800
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 1908
$ws.Range("N12").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O12").Value = 'Provincia de Diguillín'
$ws.Range("P12").Value = 1908
$ws.Range("Q12").Value = 1

$ws.Range("D13").Value = 44166
$ws.Range("J13").Value = 240
$ws.Range("K13").Value = 600
$ws.Range("L13").Value = 700
$ws.Range("M13").Value = 633
$ws.Range("N13").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O13").Value = 'Provincia de Diguillín'
$ws.Range("P13").Value = 633
$ws.Range("Q13").Value = 1

$ws.Range("D14").Value = 44273
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 1800
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 1914
$ws.Range("N14").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O14").Value = 'Provincia de Diguillín'
$ws.Range("P14").Value = 1914
$ws.Range("Q14").Value = 1

$ws.Range("D15").Value = 44160
$ws.Range("J15").Value = 190
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = 1395
$ws.Range("N15").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O15").Value = 'Provincia de Diguillín'
$ws.Range("P15").Value = 930
$ws.Range("Q15").Value = 1.5

$ws.Range("D16").Value = 44211
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 1883
$ws.Range("N16").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O16").Value = 'Provincia de Diguillín'
$ws.Range("P16").Value = 1883
$ws.Range("Q16").Value = 1
